# "Generate Report for Handback" — localization-status.xlsx
#
# The handback run for de-de completed (in sync with en-US). This:
#  - flips the Overview "Status" cells for both languages from
#    "Ready for handoff" to "Handed back: in sync with en-US"
#  - fills in the "Latest Target File" / "Latest Handback File" hyperlinked
#    columns on both the zh-cn and de-de detail sheets (previously blank)
#  - stamps a real "Latest Handback DateTime" for de-de (the language that
#    was actually handed back) and backfills the zh-cn placeholder from the
#    zero-date sentinel to the handoff-adjacent timestamp
#  - widens a few columns that now hold long file names / datetimes

$wb = $excel.ActiveWorkbook

$commitSha = "9e158e76f003be25f4bc28ee62e1cf543e208665"
$repoBase  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/"

$docA = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.md"
$docB = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.md"
$urlA = $repoBase + $docA
$urlB = $repoBase + $docB

# ColumnWidth (character units) values chosen so the exported OOXML <col>
# width lands on the intended target after the host's internal px rounding.
$wideWidth   = 29.166666666666668   # -> exported width ~29.98 (was ~17.22)
$fortyWidth  = 39.166666666666664   # -> exported width 40    (Latest Target/Handback File)

# ---------------------------------------------------------------------
# 1) Overview sheet: Status columns (zh-cn / de-de) for both doc rows
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Cells.Item(2, 5).Value = "Handed back: in sync with en-US"   # E2 (zh-cn status, row eb4d2096)
$ov.Cells.Item(2, 6).Value = "Handed back: in sync with en-US"   # F2 (de-de status, row eb4d2096)
$ov.Cells.Item(3, 5).Value = "Handed back: in sync with en-US"   # E3 (zh-cn status, row fe0a9c1a)
$ov.Cells.Item(3, 6).Value = "Handed back: in sync with en-US"   # F3 (de-de status, row fe0a9c1a)

$ov.Columns.Item(5).ColumnWidth = $wideWidth
$ov.Columns.Item(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# 2) zh-cn detail sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Columns.Item(3).ColumnWidth  = $wideWidth    # C  Status
$zh.Columns.Item(9).ColumnWidth  = $fortyWidth   # I  Latest Target File
$zh.Columns.Item(10).ColumnWidth = $fortyWidth   # J  Latest Handback File

# Row 2 (eb4d2096...)
$zh.Hyperlinks.Add($zh.Cells.Item(2, 9), $urlA, "", "", $docA)
$zh.Cells.Item(2, 10).Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.5387253624e2db7618f5b0610dc21bb69a71115a.zh-cn.xlf"

# Row 3 (fe0a9c1a...)
$zh.Hyperlinks.Add($zh.Cells.Item(3, 9), $urlB, "", "", $docB)
$zh.Cells.Item(3, 10).Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.8271afb11851a3eb9dd5f5ed7a864122870a8f1d.zh-cn.xlf"

# The zh-cn "Latest Handback DateTime" cells still share the old sentinel
# shared string ("0001-01-01 00:00:00") — back-fill it to the real
# handoff-adjacent stamp for both rows.
$zh.Cells.Item(2, 11).Value = "2016-08-28 11:09:47"
$zh.Cells.Item(3, 11).Value = "2016-08-28 11:09:47"

# ---------------------------------------------------------------------
# 3) de-de detail sheet — this is the language actually handed back
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Columns.Item(3).ColumnWidth  = $wideWidth    # C  Status
$de.Columns.Item(9).ColumnWidth  = $fortyWidth   # I  Latest Target File
$de.Columns.Item(10).ColumnWidth = $fortyWidth   # J  Latest Handback File

# Row 2 (eb4d2096...)
$de.Hyperlinks.Add($de.Cells.Item(2, 9), $urlA, "", "", $docA)
$de.Cells.Item(2, 10).Value = "eb4d2096-cf7a-4f78-a4ec-d67885d88521.5387253624e2db7618f5b0610dc21bb69a71115a.de-de.xlf"
$de.Cells.Item(2, 11).Value = "2016-08-28 11:09:54"

# Row 3 (fe0a9c1a...)
$de.Hyperlinks.Add($de.Cells.Item(3, 9), $urlB, "", "", $docB)
$de.Cells.Item(3, 10).Value = "fe0a9c1a-c85e-41cb-9bb5-42f90f64a307.8271afb11851a3eb9dd5f5ed7a864122870a8f1d.de-de.xlf"
$de.Cells.Item(3, 11).Value = "2016-08-28 11:09:54"
